$wb = $excel.ActiveWorkbook
$wsMeans = $wb.Worksheets.Item("Means")
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# --- Header label updates (shared strings) ---
$wsMeans.Range("B1").Value = "National Average"
$wsMeans.Range("C1").Value = "State Average"
$wsSD.Range("B1").Value = "National Average SD"
$wsSD.Range("C1").Value = "State Average SD"

# --- Means sheet numeric updates ---
$wsMeans.Range("B2").Value = 72
$wsMeans.Range("C2").Value = 74
$wsMeans.Range("B3").Value = 13
$wsMeans.Range("C3").Value = 12
$wsMeans.Range("E3").Value = 2.2
$wsMeans.Range("F3").Value = 1.8
$wsMeans.Range("G3").Value = 2.1
$wsMeans.Range("B4").Value = 15
$wsMeans.Range("C4").Value = 14
$wsMeans.Range("E4").Value = 6.2
$wsMeans.Range("F4").Value = 6.5
$wsMeans.Range("G4").Value = 7.1
$wsMeans.Range("B5").Value = 18
$wsMeans.Range("C5").Value = 39
$wsMeans.Range("E5").Value = 43
$wsMeans.Range("G5").Value = 40
$wsMeans.Range("B6").Value = 71
$wsMeans.Range("C6").Value = 69
$wsMeans.Range("E6").Value = 77
$wsMeans.Range("F6").Value = 70
$wsMeans.Range("G6").Value = 61
$wsMeans.Range("B7").Value = 7.3
$wsMeans.Range("C7").Value = 8.2
$wsMeans.Range("E7").Value = 4.4
$wsMeans.Range("F7").Value = 3.1
$wsMeans.Range("G7").Value = 6
$wsMeans.Range("B8").Value = 5.8
$wsMeans.Range("C8").Value = 6.2
$wsMeans.Range("E8").Value = 3.4
$wsMeans.Range("F8").Value = 3.9
$wsMeans.Range("B9").Value = 29
$wsMeans.Range("C9").Value = 31
$wsMeans.Range("B10").Value = 0.37
$wsMeans.Range("C10").Value = 0.36

# --- Standard Deviations sheet numeric updates ---
$wsSD.Range("B2").Value = 27
$wsSD.Range("C2").Value = 21
$wsSD.Range("E2").Value = 5.2
$wsSD.Range("F2").Value = 6.3
$wsSD.Range("G2").Value = 8.3
$wsSD.Range("B3").Value = 23
$wsSD.Range("C3").Value = 17
$wsSD.Range("E3").Value = 1.7
$wsSD.Range("G3").Value = 2.9
$wsSD.Range("B4").Value = 16
$wsSD.Range("C4").Value = 12
$wsSD.Range("E4").Value = 3.9
$wsSD.Range("F4").Value = 6
$wsSD.Range("G4").Value = 8.3
$wsSD.Range("B5").Value = 22
$wsSD.Range("C5").Value = 30
$wsSD.Range("E5").Value = 21
$wsSD.Range("F5").Value = 21
$wsSD.Range("G5").Value = 22
$wsSD.Range("B6").Value = 37
$wsSD.Range("C6").Value = 36
$wsSD.Range("E6").Value = 19
$wsSD.Range("F6").Value = 16
$wsSD.Range("G6").Value = 16
$wsSD.Range("B7").Value = 8.7
$wsSD.Range("C7").Value = 9.7
$wsSD.Range("E7").Value = 4.8
$wsSD.Range("F7").Value = 3.7
$wsSD.Range("B8").Value = 7.8
$wsSD.Range("C8").Value = 8.2
$wsSD.Range("E8").Value = 3.6
$wsSD.Range("G8").Value = 8.9
$wsSD.Range("B9").Value = 10
$wsSD.Range("C9").Value = 14
$wsSD.Range("C10").Value = 0.079
$wsSD.Range("E10").Value = 0.03
$wsSD.Range("F10").Value = 0.03
$wsSD.Range("G10").Value = 0.033
